# Auto-generated edit script: applies numeric value corrections to multiple sheets
# as produced by a scheduled recompute/runner over the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1519.3334
$ws.Range("J17").Value = 1519.3334
$ws.Range("L17").Value = 4558.0002
$ws.Range("N17").Value = -4894.0002
$ws.Range("H19").Value = 2802.5
$ws.Range("I19").Value = 2464.1428
$ws.Range("J19").Value = 3065.6667
$ws.Range("K19").Value = 2464.1428
$ws.Range("L19").Value = 3065.6667
$ws.Range("M19").Value = -2289.1428
$ws.Range("N19").Value = -3415.6667
$ws.Range("H69").Value = 19785.428
$ws.Range("I69").Value = 18498
$ws.Range("K69").Value = 55494
$ws.Range("M69").Value = -54620
$ws.Range("H70").Value = 9529.643
$ws.Range("I70").Value = 2739.2
$ws.Range("J70").Value = 13302.111
$ws.Range("K70").Value = 8217.599999999999
$ws.Range("L70").Value = 39906.333
$ws.Range("M70").Value = -7947.599999999999
$ws.Range("N70").Value = -40446.333
$ws.Range("H72").Value = 19785.428
$ws.Range("I72").Value = 18498
$ws.Range("K72").Value = 166482
$ws.Range("M72").Value = -162114
$ws.Range("H73").Value = 9529.643
$ws.Range("I73").Value = 2739.2
$ws.Range("J73").Value = 13302.111
$ws.Range("K73").Value = 8217.599999999999
$ws.Range("L73").Value = 39906.333
$ws.Range("M73").Value = -7281.599999999999
$ws.Range("N73").Value = -41778.333
$ws.Range("H106").Value = 9042.111000000001
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 1375.9286
$ws.Range("I107").Value = 1119.4546
$ws.Range("K107").Value = 1119.4546
$ws.Range("M107").Value = 800.5454
$ws.Range("H125").Value = 1499
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H138").Value = 2640.6792
$ws.Range("J138").Value = 3177.0513
$ws.Range("L138").Value = 9531.153900000001
$ws.Range("N138").Value = -19811.1539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 14554
$ws.Range("I46").Value = 9586.6
$ws.Range("J46").Value = 19521.4
$ws.Range("K46").Value = 9586.6
$ws.Range("L46").Value = 19521.4
$ws.Range("M46").Value = -9267.6
$ws.Range("N46").Value = -20159.4
$ws.Range("H74").Value = 1646.5588
$ws.Range("I74").Value = 1647.9697
$ws.Range("K74").Value = 1647.9697
$ws.Range("M74").Value = -773.9697000000001
$ws.Range("H77").Value = 1646.5588
$ws.Range("I77").Value = 1647.9697
$ws.Range("K77").Value = 8239.8485
$ws.Range("M77").Value = -3871.8485
$ws.Range("H132").Value = 2039.7567
$ws.Range("I132").Value = 1943.3704
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 5830.1112
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = -3300.1112
$ws.Range("N132").Value = -11960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 18259.834
$ws.Range("J106").Value = 18259.834
$ws.Range("L106").Value = 18259.834
$ws.Range("N106").Value = -20783.834
$ws.Range("H134").Value = 2533.147
$ws.Range("I134").Value = 2533.147
$ws.Range("K134").Value = 7599.441
$ws.Range("M134").Value = -5064.441

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2580200
$ws.Range("I6").Value = 2580200
$ws.Range("K6").Value = 2580200
$ws.Range("M6").Value = -2580087
$ws.Range("H31").Value = 4814.972
$ws.Range("I31").Value = 7682.636
$ws.Range("J31").Value = 3553.2
$ws.Range("K31").Value = 7682.636
$ws.Range("L31").Value = 3553.2
$ws.Range("M31").Value = -7387.636
$ws.Range("N31").Value = -4143.2
$ws.Range("H34").Value = 4814.972
$ws.Range("I34").Value = 7682.636
$ws.Range("J34").Value = 3553.2
$ws.Range("K34").Value = 7682.636
$ws.Range("L34").Value = 3553.2
$ws.Range("M34").Value = -7480.636
$ws.Range("N34").Value = -3957.2
$ws.Range("H60").Value = 6000
$ws.Range("I60").Value = 6000
$ws.Range("K60").Value = 6000
$ws.Range("M60").Value = -5489
$ws.Range("H86").Value = 14707.286
$ws.Range("I86").Value = 25481.4
$ws.Range("J86").Value = 8721.666999999999
$ws.Range("K86").Value = 25481.4
$ws.Range("L86").Value = 8721.666999999999
$ws.Range("M86").Value = -24358.4
$ws.Range("N86").Value = -10967.667
$ws.Range("H89").Value = 14707.286
$ws.Range("I89").Value = 25481.4
$ws.Range("J89").Value = 8721.666999999999
$ws.Range("K89").Value = 127407
$ws.Range("L89").Value = 43608.335
$ws.Range("M89").Value = -121791
$ws.Range("N89").Value = -54840.335
$ws.Range("H96").Value = 21116.715
$ws.Range("J96").Value = 21116.715
$ws.Range("L96").Value = 21116.715
$ws.Range("N96").Value = -26608.715
$ws.Range("H134").Value = 1203.3334
$ws.Range("I134").Value = 971.8182
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 2915.4546
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -380.4546
$ws.Range("N134").Value = -16320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69927
$ws.Range("J37").Value = 69927
$ws.Range("L37").Value = 209781
$ws.Range("N37").Value = -210005
$ws.Range("H82").Value = 613
$ws.Range("I82").Value = 613
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1839
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("M82").Value = -1433
$ws.Range("H85").Value = 613
$ws.Range("I85").Value = 613
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1839
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("M85").Value = -435
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3431.2222
$ws.Range("I80").Value = 3586
$ws.Range("K80").Value = 3586
$ws.Range("M80").Value = -2588
$ws.Range("H83").Value = 3431.2222
$ws.Range("I83").Value = 3586
$ws.Range("K83").Value = 17930
$ws.Range("M83").Value = -12938
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920
$ws.Range("H132").Value = 3863.9038
$ws.Range("I132").Value = 4106.05
$ws.Range("K132").Value = 12318.15
$ws.Range("M132").Value = -9788.150000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 442.65625
$ws.Range("I55").Value = 514.53845
$ws.Range("K55").Value = 514.53845
$ws.Range("M55").Value = -341.53845
$ws.Range("H61").Value = 2864.5
$ws.Range("I61").Value = 2507.6667
$ws.Range("K61").Value = 2507.6667
$ws.Range("M61").Value = -2305.6667
$ws.Range("H82").Value = 2849
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2849
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2849
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3571
$ws.Range("H85").Value = 2849
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2849
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2849
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5345
$ws.Range("H93").Value = 1933.2858
$ws.Range("I93").Value = 1935.8462
$ws.Range("J93").Value = 1900
$ws.Range("K93").Value = 1935.8462
$ws.Range("L93").Value = 1900
$ws.Range("M93").Value = -687.8462
$ws.Range("N93").Value = -4396
$ws.Range("H100").Value = 3716.389
$ws.Range("I100").Value = 2221.7778
$ws.Range("K100").Value = 2221.7778
$ws.Range("M100").Value = -1680.7778
$ws.Range("H113").Value = 2864.5
$ws.Range("I113").Value = 2507.6667
$ws.Range("K113").Value = 2507.6667
$ws.Range("M113").Value = -337.6667000000002
$ws.Range("H132").Value = 2510.739
$ws.Range("I132").Value = 2260.3125
$ws.Range("K132").Value = 6780.9375
$ws.Range("M132").Value = -4250.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 24911.4
$ws.Range("J74").Value = 24911.4
$ws.Range("L74").Value = 24911.4
$ws.Range("N74").Value = -26783.4
$ws.Range("H77").Value = 24911.4
$ws.Range("J77").Value = 24911.4
$ws.Range("L77").Value = 74734.20000000001
$ws.Range("N77").Value = -84094.20000000001
$ws.Range("H107").Value = 3834.9375
$ws.Range("I107").Value = 3375.6428
$ws.Range("K107").Value = 10126.9284
$ws.Range("M107").Value = -8206.928400000001
$ws.Range("H117").Value = 53502.668
$ws.Range("J117").Value = 53502.668
$ws.Range("L117").Value = 53502.668
$ws.Range("N117").Value = -62680.668
$ws.Range("H122").Value = 2772.875
$ws.Range("I122").Value = 2531.3333
$ws.Range("K122").Value = 7593.999899999999
$ws.Range("M122").Value = -5143.999899999999
$ws.Range("H126").Value = 2257.7778
$ws.Range("I126").Value = 2257.7778
$ws.Range("K126").Value = 6773.3334
$ws.Range("M126").Value = -4303.3334
$ws.Range("H132").Value = 4394.9688
$ws.Range("I132").Value = 3777.2083
$ws.Range("K132").Value = 11331.6249
$ws.Range("M132").Value = -8801.624899999999
$ws.Range("H136").Value = 1343.6923
$ws.Range("I136").Value = 1158.1316
$ws.Range("K136").Value = 3474.3948
$ws.Range("M136").Value = -924.3948
